$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits ---
$ws.Range("B2").Value = "Закрыть клапан насоса"
$ws.Range("C2").Value = ""

# --- Row 3 edits ---
$ws.Range("B3").Value = "Открыть клапан"
$ws.Range("C3").Value = "C_2"
$ws.Range("D3").Value = ""

# --- New Row 4 (values, forced to text where numeric-looking) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Значение на РРГ с паузой"
$ws.Range("C4").Value = "N_2"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.0"
$ws.Range("D4").Style = $ws.Range("D2").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0030:02"
$ws.Range("E4").Style = $ws.Range("E2").Style
$ws.Range("F4").Value = "kuda ti"

# --- New Row 5 ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Стабилизировать давление"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "1e-4"
$ws.Range("C5").Style = $ws.Range("C3").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0000:00"
$ws.Range("D5").Style = $ws.Range("D3").Style
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "ewrr"

# Copy the numeric "index" column style (bold, bordered, centered) from row 3
# down to the two new rows' A cells so they match the rest of column A.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
